# Generate Report for Handoff
# Updates the localization-status report: the zh-cn / de-de items are now
# ready for handoff (status text) and the handoff timestamps advance.
# Also narrows the "Status" columns to fit the new, shorter status text.

$wb = $excel.ActiveWorkbook

$newStatus    = "Ready for handoff"
$statusWidth  = 16.33   # renders to the narrow "Ready for handoff" column width

# ---- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus            # zh-cn status
$overview.Range("F2").Value = $newStatus            # de-de status
$overview.Range("G2").Value = "2016-08-27 08:58:53" # Latest HO Xliff Generate Date

$overview.Columns.Item(5).ColumnWidth = $statusWidth
$overview.Columns.Item(6).ColumnWidth = $statusWidth

# ---- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus                 # Status
$zhcn.Range("H2").Value = "2016-08-27 08:58:49"      # Latest Handoff Datetime

$zhcn.Columns.Item(3).ColumnWidth = $statusWidth

# ---- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus                 # Status
$dede.Range("H2").Value = "2016-08-27 08:58:53"      # Latest Handoff Datetime

$dede.Columns.Item(3).ColumnWidth = $statusWidth
